$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived NATMI values for the Vcam1-Itgb1 ligand-receptor pair sheet.
# Underlying source values changed (ECs ligand/receptor expression), and all dependent
# specificity / edge-weight columns below were recomputed accordingly.
$values = @{
    "G2" = 3.037522333333333
    "H2" = 9.112567
    "I2" = 0.1153015356242242
    "J2" = 0.1153015356242242
    "M2" = 77.08952333333333
    "N2" = 231.26857
    "O2" = 0.2403816673726824
    "P2" = 0.2403816673726824
    "Q2" = 234.1611487910211
    "R2" = 2107.45033911919
    "S2" = 0.02771637538398175
    "T2" = 0.02771637538398175
    "G3" = 3.037522333333333
    "H3" = 9.112567
    "I3" = 0.1153015356242242
    "J3" = 0.1153015356242242
    "O3" = 0.3167483425780597
    "P3" = 0.3167483425780597
    "Q3" = 308.5516320208338
    "R3" = 2776.964688187504
    "S3" = 0.03652157030567812
    "T3" = 0.03652157030567812
    "G4" = 3.037522333333333
    "H4" = 9.112567
    "I4" = 0.1153015356242242
    "J4" = 0.1153015356242242
    "O4" = 0.4428699900492579
    "P4" = 0.4428699900492579
    "Q4" = 431.4095445316285
    "R4" = 3882.685900784656
    "S4" = 0.05106358993456433
    "T4" = 0.05106358993456433
    "I5" = 0.325995654495798
    "J5" = 0.325995654495798
    "M5" = 77.08952333333333
    "N5" = 231.26857
    "O5" = 0.2403816673726824
    "P5" = 0.2403816673726824
    "Q5" = 662.0511734240877
    "R5" = 5958.460560816789
    "S5" = 0.07836337898394882
    "T5" = 0.07836337898394882
    "I6" = 0.325995654495798
    "J6" = 0.325995654495798
    "O6" = 0.3167483425780597
    "P6" = 0.3167483425780597
    "S6" = 0.1032585832491938
    "T6" = 0.1032585832491938
    "I7" = 0.325995654495798
    "J7" = 0.325995654495798
    "O7" = 0.4428699900492579
    "P7" = 0.4428699900492579
    "S7" = 0.1443736922626554
    "T7" = 0.1443736922626554
    "I8" = 0.5587028098799778
    "J8" = 0.5587028098799777
    "M8" = 77.08952333333333
    "N8" = 231.26857
    "O8" = 0.2403816673726824
    "P8" = 0.2403816673726824
    "Q8" = 1134.646568981005
    "R8" = 10211.81912082904
    "S8" = 0.1343019130047518
    "T8" = 0.1343019130047518
    "I9" = 0.5587028098799778
    "J9" = 0.5587028098799777
    "O9" = 0.3167483425780597
    "P9" = 0.3167483425780597
    "R9" = 13456.00443071527
    "S9" = 0.1769681890231878
    "T9" = 0.1769681890231877
    "I10" = 0.5587028098799778
    "J10" = 0.5587028098799777
    "O10" = 0.4428699900492579
    "P10" = 0.4428699900492579
    "R10" = 18813.8649750473
    "S10" = 0.2474327078520382
    "T10" = 0.2474327078520382
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}

